# Auto-generated edit script: expand LR-pairs table from 4x3 to 4x4 sending/target cluster matrix
# with updated TPM-derived values (commit: "update scripts wuth new tpm")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
    @{A="ECs"; D="ECs"; E=3; F=1; G=1.758258666666666; H=5.274775999999999; I=0.1132279568112417; J=0.1132279568112417; K=1; L=0.3333333333333333; M=0.009559333333333333; N=0.028678; O=0.001091248180460288; P=0.001091248180460288; Q=0.01680778068088888; R=0.151270026128; S=0.0001235598018475035; T=0.0001235598018475035},
    @{A="ECs"; D="FAPs"; E=3; F=1; G=1.758258666666666; H=5.274775999999999; I=0.1132279568112417; J=0.1132279568112417; K=3; L=1; M=8.501068666666667; N=25.503206; O=0.9704417024689271; P=0.9704417024689271; Q=14.94707765909511; R=134.523698931856; S=0.1098811311749795; T=0.1098811311749795},
    @{A="ECs"; D="MuSCs"; E=3; F=1; G=1.758258666666666; H=5.274775999999999; I=0.1132279568112417; J=0.1132279568112417; K=3; L=1; M=0.1915436666666667; N=0.574631; O=0.02186571703696477; P=0.02186571703696477; Q=0.3367833119617777; R=3.031049807656; S=0.002475810464308278; T=0.002475810464308278},
    @{A="ECs"; D="Resolving-Mac"; E=3; F=1; G=1.758258666666666; H=5.274775999999999; I=0.1132279568112417; J=0.1132279568112417; K=1; L=0.3333333333333333; M=0.05782766666666667; N=0.173483; O=0.006601332313647817; P=0.006601332313647817; Q=0.1016759960897778; R=0.9150839648079998; S=0.0007474553701063691; T=0.0007474553701063691},
    @{A="FAPs"; D="ECs"; E=3; F=1; G=2.690388333333333; H=8.071165; I=0.1732550390834427; J=0.1732550390834427; K=1; L=0.3333333333333333; M=0.009559333333333333; N=0.028678; O=0.001091248180460288; P=0.001091248180460288; Q=0.02571831887444444; R=0.23146486987; S=0.0001890642461553829; T=0.0001890642461553829},
    @{A="FAPs"; D="FAPs"; E=3; F=1; G=2.690388333333333; H=8.071165; I=0.1732550390834427; J=0.1732550390834427; K=3; L=1; M=8.501068666666667; N=25.503206; O=0.9704417024689271; P=0.9704417024689271; Q=22.87117596166556; R=205.84058365499; S=0.1681339150894566; T=0.1681339150894566},
    @{A="FAPs"; D="MuSCs"; E=3; F=1; G=2.690388333333333; H=8.071165; I=0.1732550390834427; J=0.1732550390834427; K=3; L=1; M=0.1915436666666667; N=0.574631; O=0.02186571703696477; P=0.02186571703696477; Q=0.5153268461238889; R=4.637941615115; S=0.00378834565982683; T=0.00378834565982683},
    @{A="FAPs"; D="Resolving-Mac"; E=3; F=1; G=2.690388333333333; H=8.071165; I=0.1732550390834427; J=0.1732550390834427; K=1; L=0.3333333333333333; M=0.05782766666666667; N=0.173483; O=0.006601332313647817; P=0.006601332313647817; Q=0.1555788797438889; R=1.400209917695; S=0.001143714088003846; T=0.001143714088003846},
    @{A="MuSCs"; D="ECs"; E=3; F=1; G=9.938311666666667; H=29.814935; I=0.6400052196548212; J=0.640005219654821; K=1; L=0.3333333333333333; M=0.009559333333333333; N=0.028678; O=0.001091248180460288; P=0.001091248180460288; Q=0.09500363399222223; R=0.85503270593; S=0.0006984045314334103; T=0.0006984045314334102},
    @{A="MuSCs"; D="FAPs"; E=3; F=1; G=9.938311666666667; H=29.814935; I=0.6400052196548212; J=0.640005219654821; K=3; L=1; M=8.501068666666667; N=25.503206; O=0.9704417024689271; P=0.9704417024689271; Q=84.48626990906779; R=760.37642918161; S=0.6210877549508244; T=0.6210877549508242},
    @{A="MuSCs"; D="MuSCs"; E=3; F=1; G=9.938311666666667; H=29.814935; I=0.6400052196548212; J=0.640005219654821; K=3; L=1; M=0.1915436666666667; N=0.574631; O=0.02186571703696477; P=0.02186571703696477; Q=1.903620657109445; R=17.132585913985; S=0.0139941730351528; T=0.0139941730351528},
    @{A="MuSCs"; D="Resolving-Mac"; E=3; F=1; G=9.938311666666667; H=29.814935; I=0.6400052196548212; J=0.640005219654821; K=1; L=0.3333333333333333; M=0.05782766666666667; N=0.173483; O=0.006601332313647817; P=0.006601332313647817; Q=0.5747093742894445; R=5.172384368605; S=0.004224887137410639; T=0.004224887137410639},
    @{A="Resolving-Mac"; D="ECs"; E=3; F=1; G=1.141526666666667; H=3.42458; I=0.0735117844504946; J=0.07351178445049458; K=1; L=0.3333333333333333; M=0.009559333333333333; N=0.028678; O=0.001091248180460288; P=0.001091248180460288; Q=0.01091223391555555; R=0.09821010523999998; S=0.0000802196010239911; T=0.00008021960102399109},
    @{A="Resolving-Mac"; D="FAPs"; E=3; F=1; G=1.141526666666667; H=3.42458; I=0.0735117844504946; J=0.07351178445049458; K=3; L=1; M=8.501068666666667; N=25.503206; O=0.9704417024689271; P=0.9704417024689271; Q=9.704196578164444; R=87.33776920348; S=0.07133890125366678; T=0.07133890125366676},
    @{A="Resolving-Mac"; D="MuSCs"; E=3; F=1; G=1.141526666666667; H=3.42458; I=0.0735117844504946; J=0.07351178445049458; K=3; L=1; M=0.1915436666666667; N=0.574631; O=0.02186571703696477; P=0.02186571703696477; Q=0.2186522033311111; R=1.96786982998; S=0.001607387877676861; T=0.001607387877676861},
    @{A="Resolving-Mac"; D="Resolving-Mac"; E=3; F=1; G=1.141526666666667; H=3.42458; I=0.0735117844504946; J=0.07351178445049458; K=1; L=0.3333333333333333; M=0.05782766666666667; N=0.173483; O=0.006601332313647817; P=0.006601332313647817; Q=0.0660118235711111; R=0.5941064121399999; S=0.0004852757181269631; T=0.000485275718126963},
)

$r = 2
foreach ($row in $rowsData) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "Clcf1"
    $ws.Cells.Item($r, 3).Value = "Cntfr"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}

"done: wrote $($rowsData.Count) rows"
